$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 201.46666
$ws.Range("I9").Value = 210.7
$ws.Range("J9").Value = 183
$ws.Range("K9").Value = 210.7
$ws.Range("L9").Value = 183
$ws.Range("M9").Value = -41.69999999999999
$ws.Range("N9").Value = -521
# Row 40
$ws.Range("H40").Value = 2931.9285
$ws.Range("I40").Value = 2549.8333
$ws.Range("J40").Value = 3218.5
$ws.Range("K40").Value = 2549.8333
$ws.Range("L40").Value = 3218.5
$ws.Range("M40").Value = -2374.8333
$ws.Range("N40").Value = -3568.5
# Row 51
$ws.Range("H51").Value = 4000
$ws.Range("J51").Value = 4000
$ws.Range("L51").Value = 4000
$ws.Range("N51").Value = -4968
# Row 62
$ws.Range("H62").Value = 9405.154
$ws.Range("I62").Value = 6327.9
$ws.Range("K62").Value = 6327.9
$ws.Range("M62").Value = -5703.9
# Row 65
$ws.Range("H65").Value = 9405.154
$ws.Range("I65").Value = 6327.9
$ws.Range("K65").Value = 31639.5
$ws.Range("M65").Value = -28519.5
# Row 70
$ws.Range("H70").Value = 9529110
$ws.Range("I70").Value = 25004710
$ws.Range("K70").Value = 75014130
$ws.Range("M70").Value = -75013860
# Row 73
$ws.Range("H73").Value = 9529110
$ws.Range("I73").Value = 25004710
$ws.Range("K73").Value = 75014130
$ws.Range("M73").Value = -75013194
# Row 92
$ws.Range("H92").Value = 437.0476
$ws.Range("I92").Value = 398.8
$ws.Range("J92").Value = 532.6667
$ws.Range("K92").Value = 398.8
$ws.Range("L92").Value = 532.6667
$ws.Range("M92").Value = 849.2
$ws.Range("N92").Value = -3028.6667
# Row 106
$ws.Range("H106").Value = 4117754
$ws.Range("I106").Value = 5880420
$ws.Range("J106").Value = 4866.222
$ws.Range("K106").Value = 5880420
$ws.Range("L106").Value = 4866.222
$ws.Range("M106").Value = -5879789
$ws.Range("N106").Value = -6128.222
# Row 129
$ws.Range("H129").Value = 125000650
$ws.Range("I129").Value = 737.4286
$ws.Range("J129").Value = 1000000000
$ws.Range("K129").Value = 2212.2858
$ws.Range("L129").Value = 3000000000
$ws.Range("M129").Value = 2787.7142
$ws.Range("N129").Value = -3000010000
# Row 137
$ws.Range("H137").Value = 759968.2
$ws.Range("I137").Value = 1503253.2
$ws.Range("J137").Value = 16683.143
$ws.Range("K137").Value = 4509759.6
$ws.Range("L137").Value = 50049.429
$ws.Range("M137").Value = -4507209.6
$ws.Range("N137").Value = -55149.429
# Row 138
$ws.Range("H138").Value = 4807.0923
$ws.Range("I138").Value = 514.5
$ws.Range("J138").Value = 6139.276
$ws.Range("K138").Value = 1543.5
$ws.Range("L138").Value = 18417.828
$ws.Range("M138").Value = 3596.5
$ws.Range("N138").Value = -28697.828
# Row 141
$ws.Range("H141").Value = 5203.125
$ws.Range("I141").Value = 4769.875
$ws.Range("J141").Value = 6069.625
$ws.Range("K141").Value = 14309.625
$ws.Range("L141").Value = 18208.875
$ws.Range("M141").Value = -9129.625
$ws.Range("N141").Value = -28568.875

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2598.9302
$ws.Range("I32").Value = 2554.9756
$ws.Range("K32").Value = 2554.9756
$ws.Range("M32").Value = -2267.9756
# Row 45
$ws.Range("H45").Value = 14161.6
$ws.Range("I45").Value = 20349.75
$ws.Range("J45").Value = 7089.4287
$ws.Range("K45").Value = 20349.75
$ws.Range("L45").Value = 7089.4287
$ws.Range("M45").Value = -19972.75
$ws.Range("N45").Value = -7843.4287
# Row 46
$ws.Range("H46").Value = 8614.833000000001
$ws.Range("I46").Value = 4998
$ws.Range("J46").Value = 9338.200000000001
$ws.Range("K46").Value = 4998
$ws.Range("L46").Value = 9338.200000000001
$ws.Range("M46").Value = -4679
$ws.Range("N46").Value = -9976.200000000001
# Row 74
$ws.Range("H74").Value = 4252.3477
$ws.Range("I74").Value = 2150.2354
$ws.Range("J74").Value = 10208.333
$ws.Range("K74").Value = 2150.2354
$ws.Range("L74").Value = 10208.333
$ws.Range("M74").Value = -1276.2354
$ws.Range("N74").Value = -11956.333
# Row 77
$ws.Range("H77").Value = 4252.3477
$ws.Range("I77").Value = 2150.2354
$ws.Range("J77").Value = 10208.333
$ws.Range("K77").Value = 10751.177
$ws.Range("L77").Value = 51041.665
$ws.Range("M77").Value = -6383.177
$ws.Range("N77").Value = -59777.665
# Row 125
$ws.Range("H125").Value = 173000
$ws.Range("J125").Value = 173000
$ws.Range("L125").Value = 173000
$ws.Range("N125").Value = -182840

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 4395.08
$ws.Range("I86").Value = 6395.643
$ws.Range("J86").Value = 1848.909
$ws.Range("K86").Value = 6395.643
$ws.Range("L86").Value = 1848.909
$ws.Range("M86").Value = -5272.643
$ws.Range("N86").Value = -4094.909
# Row 89
$ws.Range("H89").Value = 4395.08
$ws.Range("I89").Value = 6395.643
$ws.Range("J89").Value = 1848.909
$ws.Range("K89").Value = 31978.215
$ws.Range("L89").Value = 9244.545
$ws.Range("M89").Value = -26362.215
$ws.Range("N89").Value = -20476.545
# Row 99
$ws.Range("H99").Value = 17541.938
$ws.Range("I99").Value = 19547.928
$ws.Range("K99").Value = 19547.928
$ws.Range("M99").Value = -18049.928
# Row 107
$ws.Range("H107").Value = 1913.5385
$ws.Range("I107").Value = 2009.1
$ws.Range("K107").Value = 2009.1
$ws.Range("M107").Value = -89.09999999999991
# Row 122
$ws.Range("H122").Value = 113777.5
$ws.Range("J122").Value = 113777.5
$ws.Range("L122").Value = 113777.5
$ws.Range("N122").Value = -123577.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 7
$ws.Range("H7").Value = 2505333.8
$ws.Range("I7").Value = 5002
$ws.Range("J7").Value = 3005400
$ws.Range("K7").Value = 5002
$ws.Range("L7").Value = 3005400
$ws.Range("M7").Value = -4890
$ws.Range("N7").Value = -3005624
# Row 8
$ws.Range("H8").Value = 2505333.8
$ws.Range("I8").Value = 5002
$ws.Range("J8").Value = 3005400
$ws.Range("K8").Value = 5002
$ws.Range("L8").Value = 3005400
$ws.Range("M8").Value = -4863
$ws.Range("N8").Value = -3005678
# Row 70
$ws.Range("H70").Value = 3412.6667
$ws.Range("I70").Value = 3168.1365
$ws.Range("K70").Value = 3168.1365
$ws.Range("M70").Value = -2898.1365
# Row 73
$ws.Range("H73").Value = 3412.6667
$ws.Range("I73").Value = 3168.1365
$ws.Range("K73").Value = 3168.1365
$ws.Range("M73").Value = -2232.1365
# Row 86
$ws.Range("H86").Value = 32333
$ws.Range("J86").Value = 32333
$ws.Range("L86").Value = 32333
$ws.Range("N86").Value = -34705
# Row 89
$ws.Range("H89").Value = 32333
$ws.Range("J89").Value = 32333
$ws.Range("L89").Value = 96999
$ws.Range("N89").Value = -108855
# Row 97
$ws.Range("H97").Value = 9840.5
$ws.Range("I97").Value = 12195.833
$ws.Range("J97").Value = 4541
$ws.Range("K97").Value = 12195.833
$ws.Range("L97").Value = 4541
$ws.Range("M97").Value = -11699.833
$ws.Range("N97").Value = -5533
# Row 141
$ws.Range("H141").Value = 100428.5
$ws.Range("J141").Value = 100428.5
$ws.Range("L141").Value = 100428.5
$ws.Range("N141").Value = -110788.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 3468.2
$ws.Range("I16").Value = 2985
$ws.Range("J16").Value = 5401
$ws.Range("K16").Value = 2985
$ws.Range("L16").Value = 5401
$ws.Range("M16").Value = -2815
$ws.Range("N16").Value = -5741
# Row 55
$ws.Range("H55").Value = 1921.625
$ws.Range("I55").Value = 1562.5
$ws.Range("J55").Value = 2999
$ws.Range("K55").Value = 1562.5
$ws.Range("L55").Value = 2999
$ws.Range("M55").Value = -1389.5
$ws.Range("N55").Value = -3345
# Row 101
$ws.Range("H101").Value = 25713.889
$ws.Range("J101").Value = 25713.889
$ws.Range("L101").Value = 25713.889
$ws.Range("N101").Value = -32203.889
# Row 122
$ws.Range("H122").Value = 2853.2144
$ws.Range("I122").Value = 2120.2415
$ws.Range("K122").Value = 6360.7245
$ws.Range("M122").Value = -3910.7245

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 3
$ws.Range("H3").Value = 70316
$ws.Range("I3").Value = 64379.4
$ws.Range("K3").Value = 64379.4
$ws.Range("M3").Value = -64265.4
# Row 6
$ws.Range("H6").Value = 4001.6667
$ws.Range("I6").Value = 4502.5
$ws.Range("J6").Value = 3000
$ws.Range("K6").Value = 4502.5
$ws.Range("L6").Value = 3000
$ws.Range("M6").Value = -4387.5
$ws.Range("N6").Value = -3230
# Row 122
$ws.Range("H122").Value = 4125.2915
$ws.Range("I122").Value = 3134.6
$ws.Range("K122").Value = 9403.799999999999
$ws.Range("M122").Value = -6953.799999999999
